# fix: precedence parsing accepts function calls
#
# This LL(1) parse table lives on Sheet1: row headers (column A) name a
# grammar production, column headers (row 1) name lookahead terminals, and
# each interior cell holds the "(rule, alt)" tag chosen for that
# production/terminal pair - "(0, 0)" (shown as shared string index 40)
# means "no entry / parse error".
#
# A new alternative was inserted into the grammar (rule "(4, 23)"), which
# renumbers every rule that used to sit at/after that slot - "(4, 5)" ->
# "(4, 19)", "(4, 7)" -> "(4, 21)" and "(4, 9)" -> "(4, 24)" are the same
# table entries, just re-tagged after the insertion. The actual behavioural
# fix is that the "args2" row now also has an entry under the ")" column
# (previously empty), so a function call's argument list is allowed to be
# followed directly by ")" (i.e. zero/trailing args close the call) instead
# of only by ",".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "expression" row, funcName ("L") column: rule re-tagged (4, 5) -> (4, 19)
$ws.Range("L11").Value = "(4, 19)"

# "args" row, columns for every terminal that can start an expression
# ( "(", funcName, variable, "Unary -", "!", constant ): rule re-tagged
# (4, 7) -> (4, 21)
$ws.Range("I12").Value = "(4, 21)"
$ws.Range("L12").Value = "(4, 21)"
$ws.Range("M12").Value = "(4, 21)"
$ws.Range("Q12").Value = "(4, 21)"
$ws.Range("AD12").Value = "(4, 21)"
$ws.Range("AF12").Value = "(4, 21)"

# "args2" row, "," column: rule re-tagged (4, 9) -> (4, 23)
$ws.Range("N13").Value = "(4, 23)"

# "args2" row, ")" column: brand-new table entry (4, 24) - this is the
# actual precedence-parsing fix, letting a call's argument list end at ")".
$ws.Range("J13").Value = "(4, 24)"

# Reflect the author's view state: scrolled/zoomed in on the new entry.
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 11
$ws.Range("N13").Select() | Out-Null
$win.Zoom = 110
